$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.222.60'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '1.803.80'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = "'314.94"
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('D7').Value = "'0.5246"
$ws.Range('E7').Value = '  +2.67%  '
$ws.Range('D8').Value = "'0.3810"
$ws.Range('E8').Value = '  -3.06%  '
$ws.Range('D9').Value = "'0.07926"
$ws.Range('E9').Value = '  +3.23%  '
$ws.Range('D10').Value = "'41.45"
$ws.Range('E10').Value = '  -1.13%  '
$ws.Range('D11').Value = "'1.095"
$ws.Range('E11').Value = '  -0.89%  '
$ws.Range('D12').Value = "'6.344"
$ws.Range('E12').Value = '  +1.15%  '
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').Value = "'20.61"
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').Value = "'7.345"
$ws.Range('E15').Value = '  -2.11%  '
$ws.Range('D16').Value = '1.805.76'
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').Value = "'93.18"
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').Value = "'0.00001090"
$ws.Range('E18').Value = '  -1.12%  '
$ws.Range('D19').Value = "'0.06593"
$ws.Range('E19').Value = '  -0.69%  '
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('D21').Value = "'17.40"
$ws.Range('E21').Value = '  -1.74%  '
$ws.Range('D22').Value = "'5.964"
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('D23').Value = '28.276.20'
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('D24').Value = "'11.11"
$ws.Range('E24').Value = '  -0.47%  '
$ws.Range('D25').Value = "'2.236"
$ws.Range('E25').Value = '  -0.82%  '
$ws.Range('D26').Value = "'157.49"
$ws.Range('E26').Value = '  +1.09%  '
$ws.Range('D27').Value = "'20.48"
$ws.Range('E27').Value = '  -2.90%  '
$ws.Range('B28').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C28').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D28').Value = '2.009.98'
$ws.Range('E28').Value = '  -1.11%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = "'2.400"
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').Value = "'123.07"
$ws.Range('E30').Value = '  -0.78%  '
$ws.Range('D31').Value = "'0.1103"
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('D32').Value = "'1.059"
$ws.Range('E32').Value = '  -3.87%  '
$ws.Range('E33').Value = '  +0.36%  '
$ws.Range('D34').Value = "'5.561"
$ws.Range('E34').Value = '  -1.55%  '
$ws.Range('D35').Value = "'0.07207"
$ws.Range('E35').Value = '  +1.77%  '
$ws.Range('D36').Value = "'12.20"
$ws.Range('E36').Value = '  +9.11%  '
$ws.Range('D37').Value = "'0.2165"
$ws.Range('E37').Value = '  -1.98%  '
$ws.Range('D38').Value = "'0.02309"
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('D39').Value = "'8.818"
$ws.Range('E39').Value = '  +0.57%  '
$ws.Range('D40').Value = "'5.036"
$ws.Range('E40').Value = '  -2.59%  '
$ws.Range('D41').Value = "'0.6191"
$ws.Range('E41').Value = '  -0.93%  '
$ws.Range('D42').Value = "'1.164"
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('D43').Value = "'1.378"
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('D44').Value = "'0.6022"
$ws.Range('E44').Value = '  +2.76%  '
$ws.Range('D45').Value = "'13.27"
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').Value = "'3.782"
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('D47').Value = "'125.94"
$ws.Range('D48').Value = "'1.210"
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('D49').Value = "'1.928"
$ws.Range('E49').Value = '  -2.52%  '
$ws.Range('D50').Value = "'0.06830"
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('D51').Value = "'72.70"
$ws.Range('E51').Value = '  -1.79%  '
